# Weekly price-sheet update: a new weekly record for Cilantro / Vega Modelo
# de Temuco is inserted above the existing row 259, pushing the rest of the
# historical rows (old 259..325) down by one (new 260..326).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 259 - everything at/below shifts down.
$ws.Rows(259).Insert()

# Seed the new row with the same "shape" as the row now sitting below it
# (old row 259, now row 260), then overwrite the fields that differ for
# this new weekly observation (date, volume, prices).
$src = $ws.Range("A260:R260")
$dst = $ws.Range("A259:R259")
$src.Copy($dst)

$ws.Cells.Item(259, 4).Value2 = 44642   # Fecha
$ws.Cells.Item(259, 10).Value2 = 150    # Volumen
$ws.Cells.Item(259, 11).Value2 = 5000   # Precio minimo
$ws.Cells.Item(259, 12).Value2 = 5000   # Precio maximo
$ws.Cells.Item(259, 13).Value2 = 5000   # Precio promedio ponderado
$ws.Cells.Item(259, 16).Value2 = 2500   # Precio $/Kg
